# Updated symbol list on Fri Dec 30 22:53:03 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores numeric-looking values as TEXT in the
# original workbook (t="inlineStr"). Plain `.Value = "123.45"` assignment
# would let Excel auto-coerce that into a real number (dropping
# significant trailing zeros, changing formatting, etc.), so every D-column
# write below first forces the cell's number format to Text ("@") to keep
# it a string, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $value) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Price (column D) updates ---
Set-TextValue "D2"  "245.11"
Set-TextValue "D3"  "25.20"
Set-TextValue "D4"  "5.043"
Set-TextValue "D6"  "6.538"
Set-TextValue "D7"  "3.018"
Set-TextValue "D8"  "0.8146"
Set-TextValue "D9"  "0.8355"
Set-TextValue "D10" "0.1338"
Set-TextValue "D11" "0.06950"
Set-TextValue "D13" "0.02830"
Set-TextValue "D14" "0.09400"
Set-TextValue "D15" "0.001518"
Set-TextValue "D16" "0.0005960"
Set-TextValue "D17" "0.006234"
Set-TextValue "D21" "0.1320"
Set-TextValue "D22" "3.740"
Set-TextValue "D26" "0.004285"
Set-TextValue "D40" "0.03658"
Set-TextValue "D44" "0.008194"
Set-TextValue "D45" "0.00005298"

# --- Row 27: "Best in 24h" suffix appended to the symbol label ---
$ws.Range("E27").Value = "26NitroExNTXBestin24h"

# --- Rows 41/42: KickToken and BKEXToken swap places (rank stays 40 & 41,
#     only the coin's identity/link/price/label are exchanged). ---
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1365"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D42" "0.006221"
$ws.Range("E42").Value = "41KickTokenKICK"
